# Generate Report for Handoff
# Adds a new tracked file ("897accac-7eaf-4d22-9552-966043025175.md") as
# row 9 to the Overview, zh-cn and de-de sheets, mirroring the layout of
# the existing "Ready for handoff" rows (e.g. row 8 / 728a2081-...).

$wb = $excel.ActiveWorkbook

$guid = "897accac-7eaf-4d22-9552-966043025175"
$commitHandoff = "fb172380f9e0c3a1f4460b1b1e7f23ecc4b9d182"

$statusReady = "Ready for handoff"
$mdExt = ".md"
$handoffReason = "Include"
$noHandback = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A9").Value = "$guid.md"
$ws1.Range("B9").Value = $statusReady
$ws1.Range("C9").Value = $statusReady
$ws1.Range("D9").Value = "2016-38-20 16:38:17"

$ws1.Hyperlinks.Add($ws1.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/b43a910e3c573c22ae05a56b82c899301028892a/e2e/$guid.md", "", "", "$guid.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$zhXlf = "$guid.$commitHandoff.zh-cn.xlf"

$ws2.Range("A9").Value = "$guid.md"
$ws2.Range("B9").Value = $mdExt
$ws2.Range("C9").Value = $statusReady
$ws2.Range("D9").Value = $zhXlf
$ws2.Range("E9").Value = "2016-03-20 16:38:14"
$ws2.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H9").Value = $noHandback
$ws2.Range("I9").Value = $handoffReason

$ws2.Hyperlinks.Add($ws2.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/b43a910e3c573c22ae05a56b82c899301028892a/e2e/$guid.md", "", "", "$guid.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B9"), "https://github.com/OpenLocalizationTest/oltest/blob/b43a910e3c573c22ae05a56b82c899301028892a/e2e/$guid.md", "", "", $mdExt) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bd6828b27e7b1b4949edcf4b25815f465a84d160/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", "", "", $zhXlf) | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$deXlf = "$guid.$commitHandoff.de-de.xlf"

$ws3.Range("A9").Value = "$guid.md"
$ws3.Range("B9").Value = $mdExt
$ws3.Range("C9").Value = $statusReady
$ws3.Range("D9").Value = $deXlf
$ws3.Range("E9").Value = "2016-03-20 16:38:17"
$ws3.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H9").Value = $noHandback
$ws3.Range("I9").Value = $handoffReason

$ws3.Hyperlinks.Add($ws3.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/b43a910e3c573c22ae05a56b82c899301028892a/e2e/$guid.md", "", "", "$guid.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B9"), "https://github.com/OpenLocalizationTest/oltest/blob/b43a910e3c573c22ae05a56b82c899301028892a/e2e/$guid.md", "", "", $mdExt) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e0bf6a265f1635d7e75c25bd9225ccd1fb3a8dc0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", "", "", $deXlf) | Out-Null
